# Add a new "PF/1.0.4" row to the meta-sheet.
#
# The existing header/value rows (1-2) are re-entered so that they pick
# up the worksheet's column default style, matching how the new row's
# style is resolved; the brand new row (3) is explicitly reset back to
# the workbook's "Normal" style so it carries no explicit per-cell
# style override.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the existing header (row 1) and value (row 2) rows.
$ws.Range("A1:D2").Clear()

$ws.Range("A1").Value = "dev2"
$ws.Range("B1").Value = "sit2"
$ws.Range("C1").Value = "uat2"
$ws.Range("D1").Value = "prod"

$ws.Range("A2").Value = "PF/1.0.0"
$ws.Range("B2").Value = "PF/1.0.0"
$ws.Range("C2").Value = "PF/1.0.0"
$ws.Range("D2").Value = "PF/1.0.0"

# New row: PF/1.0.4 applies to every environment ("X").
$ws.Range("A3").Value = "PF/1.0.4"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
$ws.Range("A3:D3").Style = "Normal"
